$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the hidden column's string values ("Secret Stuff") with sequential numbers 1..30
for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 6).Value = $row - 1
}

# Adjust column widths: C and E get new explicit widths, F becomes narrower and hidden
$ws.Columns.Item(3).ColumnWidth = 13
$ws.Columns.Item(5).ColumnWidth = 5.666666666666667
$ws.Columns.Item(6).ColumnWidth = 8.833333333333334
$ws.Columns.Item(6).Hidden = $true

# Update the view: select the whole F column (this also resets the scrolled
# top-left cell back to the sheet's origin)
$ws.Activate() | Out-Null
$ws.Range("F1:F1048576").Select() | Out-Null
